# Update the cryptocurrency symbol list with refreshed price/volume data
# (scraped on 2023-01-25).
#
# Rows 8-18: a new "GateToken" entry is inserted at the top of that block
# (row 8), pushing BTSEToken/MXToken/.../LEO each down by one row; the slot
# that used to hold GateToken's own data (row 18, at the bottom of the
# block) now shows LEO's refreshed figures. Coin/Link for every other row
# stay put; only Price (D) / Volume(1h) (E) get refreshed values.
#
# All of these sheet cells are plain text (e.g. "300.05", "-6.74%" are
# literal strings, not numbers/percentages) -- Excel's COM layer likes to
# auto-coerce number-looking text back into a real number/percentage (and,
# once forced back to text, tags the cell with a `Text` number format
# style). To keep the cells as plain text with their original (default)
# style, each write below forces the Text number format, assigns the
# string, then resets the cell style back to Normal so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple price/volume refreshes (Coin/Link unchanged) ---

Set-TextValue $ws.Range("D2") "299.52"
Set-TextValue $ws.Range("E2") "-7.01%"

Set-TextValue $ws.Range("D3") "35.03"
Set-TextValue $ws.Range("E3") "-3.06%"

Set-TextValue $ws.Range("D4") "4.983"
Set-TextValue $ws.Range("E4") "-2.97%"

Set-TextValue $ws.Range("D5") "0.07911"
Set-TextValue $ws.Range("E5") "-1.95%"

Set-TextValue $ws.Range("D6") "1.907"
Set-TextValue $ws.Range("E6") "-11.72%"

Set-TextValue $ws.Range("E7") "-4.19%"

# --- Rows 8-18: Coin/Link/Price/Volume all shift (GateToken inserted at
#     the top of the block, everything else shifts down by one row) ---

Set-TextValue $ws.Range("B8") "GateToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D8") "4.013"
Set-TextValue $ws.Range("E8") "-2.83%"

Set-TextValue $ws.Range("B9") "BTSEToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.926"
Set-TextValue $ws.Range("E9") "4.49%"

Set-TextValue $ws.Range("B10") "MXToken"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D10") "0.9238"
Set-TextValue $ws.Range("E10") "-0.47%"

Set-TextValue $ws.Range("B11") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.1131"
Set-TextValue $ws.Range("E11") "12.15%"

Set-TextValue $ws.Range("B12") "WazirX"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D12") "0.1827"
Set-TextValue $ws.Range("E12") "-3.12%"

Set-TextValue $ws.Range("B13") "MandalaExchangeToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D13") "0.09259"
Set-TextValue $ws.Range("E13") "-0.40%"

Set-TextValue $ws.Range("B14") "BitrueCoin"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03527"
Set-TextValue $ws.Range("E14") "-1.39%"

Set-TextValue $ws.Range("B15") "BitMartToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09872"
Set-TextValue $ws.Range("E15") "-0.66%"

Set-TextValue $ws.Range("B16") "BitForexToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001396"
Set-TextValue $ws.Range("E16") "-3.11%"

Set-TextValue $ws.Range("B17") "TigerCash"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.005767"
Set-TextValue $ws.Range("E17") "2.12%"

Set-TextValue $ws.Range("B18") "LEO"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.494"
Set-TextValue $ws.Range("E18") "1.09%"

# --- Remaining simple price/volume refreshes ---

Set-TextValue $ws.Range("D19") "0.3441"
Set-TextValue $ws.Range("E19") "2.00%"

Set-TextValue $ws.Range("E20") "-1.67%"

Set-TextValue $ws.Range("D21") "5.051"
Set-TextValue $ws.Range("E21") "-0.81%"

Set-TextValue $ws.Range("D23") "0.04502"
Set-TextValue $ws.Range("E23") "-2.19%"

Set-TextValue $ws.Range("D24") "0.001215"
Set-TextValue $ws.Range("E24") "-2.22%"

Set-TextValue $ws.Range("D25") "0.004575"
Set-TextValue $ws.Range("E25") "-3.58%"

Set-TextValue $ws.Range("D26") "0.0001251"
Set-TextValue $ws.Range("E26") "-3.75%"

Set-TextValue $ws.Range("E27") "-6.78%"

Set-TextValue $ws.Range("D39") "0.01884"
Set-TextValue $ws.Range("E39") "-4.54%"

Set-TextValue $ws.Range("D40") "0.04683"
Set-TextValue $ws.Range("E40") "-6.10%"

Set-TextValue $ws.Range("D41") "0.007627"
Set-TextValue $ws.Range("E41") "-2.47%"

Set-TextValue $ws.Range("D42") "0.009562"
Set-TextValue $ws.Range("E42") "22.12%"

Set-TextValue $ws.Range("E43") "-5.63%"

Set-TextValue $ws.Range("D44") "0.002121"
Set-TextValue $ws.Range("E44") "2.70%"

Set-TextValue $ws.Range("D45") "0.01112"
Set-TextValue $ws.Range("E45") "-5.80%"

Set-TextValue $ws.Range("D46") "0.00006015"
Set-TextValue $ws.Range("E46") "-5.55%"

Set-TextValue $ws.Range("E47") "0.09%"

Set-TextValue $ws.Range("E49") "-31.32%"

Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "0.09%"

Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "0.09%"
